# Generate Report for Handoff
# Updates the localization-status workbook with the new handoff GUID / hash
# based file names and refreshed timestamps, while preserving the existing
# hyperlink targets (only the displayed text changes).

$wb = $excel.ActiveWorkbook

$oldGuid = "b140426e-ecfd-4531-a185-43ff20de721e"
$newGuid = "b05d12c2-68c6-4d8a-a2fb-4186d2d0c1ef"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/56b74f61f510f24bc003ca9b856a4a61a6b2e92a/e2e/$oldGuid.md"

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value = "$newGuid.md"

$overview.Range("B2").Hyperlinks.Delete()
$overview.Hyperlinks.Add($overview.Range("B2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md") | Out-Null

$overview.Range("G2").Value = "2016-08-31 01:01:36"

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null

$zhcn.Range("G2").Value = "$newGuid.52b91a0d8d5f0611921e78edb3f47a2cf91dfaaa.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-08-31 01:01:32"

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $hyperlinkAddress, [Type]::Missing, [Type]::Missing, "$newGuid.md") | Out-Null

$dede.Range("G2").Value = "$newGuid.52b91a0d8d5f0611921e78edb3f47a2cf91dfaaa.de-de.xlf"
$dede.Range("H2").Value = "2016-08-31 01:01:36"
